$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 becomes what used to be row 5's data (date 44991, Primera, 50, 6000, 6000, 6000, 3000)
$ws.Range("D3").Value = 44991
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("S3").Value = 3000

# Row 4 becomes what used to be row 3's quality data (Especial, 7000, 7000, 7000, 3500); date/volume stay as before
$ws.Range("L4").Value = "Especial"
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("S4").Value = 3500

# Row 5 becomes date 45008 and volume 60 (rest unchanged)
$ws.Range("D5").Value = 45008
$ws.Range("M5").Value = 60
